# Scheduled runner: refresh computed market-price / profit columns (H, I, J, K, L, M, N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with the latest pulled prices.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 81.6
$ws.Range("J9").Value = 86
$ws.Range("L9").Value = 86
$ws.Range("N9").Value = -424
$ws.Range("H18").Value = 1657.6154
$ws.Range("J18").Value = 600
$ws.Range("L18").Value = 600
$ws.Range("N18").Value = -1168
$ws.Range("H32").Value = 3199.4
$ws.Range("I32").Value = 2999
$ws.Range("K32").Value = 2999
$ws.Range("M32").Value = -2673
$ws.Range("H62").Value = 7476.8887
$ws.Range("I62").Value = 6430.6665
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 6430.6665
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -5806.6665
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7476.8887
$ws.Range("I65").Value = 6430.6665
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 32153.3325
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -29033.3325
$ws.Range("N65").Value = -46240
$ws.Range("H98").Value = 871.2857
$ws.Range("I98").Value = 871.2857
$ws.Range("K98").Value = 871.2857
$ws.Range("M98").Value = 626.7143
$ws.Range("H112").Value = 3009.8333
$ws.Range("J112").Value = 3009.8333
$ws.Range("L112").Value = 9029.499899999999
$ws.Range("N112").Value = -11245.4999
$ws.Range("H122").Value = 871.2857
$ws.Range("I122").Value = 871.2857
$ws.Range("K122").Value = 2613.8571
$ws.Range("M122").Value = -163.8571000000002
$ws.Range("H137").Value = 2742.05
$ws.Range("I137").Value = 1420.0769
$ws.Range("K137").Value = 4260.2307
$ws.Range("M137").Value = -1710.2307

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7733.3413
$ws.Range("I32").Value = 5752.0557
$ws.Range("K32").Value = 5752.0557
$ws.Range("M32").Value = -5465.0557
$ws.Range("H45").Value = 1718.0526
$ws.Range("I45").Value = 1718.0526
$ws.Range("K45").Value = 1718.0526
$ws.Range("M45").Value = -1341.0526
$ws.Range("H56").Value = 29110
$ws.Range("J56").Value = 29110
$ws.Range("L56").Value = 29110
$ws.Range("N56").Value = -30594
$ws.Range("H61").Value = 2252.5
$ws.Range("I61").Value = 2252.5
$ws.Range("K61").Value = 2252.5
$ws.Range("M61").Value = -2040.5
$ws.Range("H74").Value = 2610.0908
$ws.Range("I74").Value = 1313.7778
$ws.Range("J74").Value = 8443.5
$ws.Range("K74").Value = 1313.7778
$ws.Range("L74").Value = 8443.5
$ws.Range("M74").Value = -439.7778000000001
$ws.Range("N74").Value = -10191.5
$ws.Range("H77").Value = 2610.0908
$ws.Range("I77").Value = 1313.7778
$ws.Range("J77").Value = 8443.5
$ws.Range("K77").Value = 6568.889
$ws.Range("L77").Value = 42217.5
$ws.Range("M77").Value = -2200.889
$ws.Range("N77").Value = -50953.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H132").Value = 1394.125
$ws.Range("I132").Value = 1394.125
$ws.Range("K132").Value = 4182.375
$ws.Range("M132").Value = -1652.375
$ws.Range("H136").Value = 2252.5
$ws.Range("I136").Value = 2252.5
$ws.Range("K136").Value = 6757.5
$ws.Range("M136").Value = -4207.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 736.125
$ws.Range("I107").Value = 725.26666
$ws.Range("J107").Value = 899
$ws.Range("K107").Value = 725.26666
$ws.Range("L107").Value = 899
$ws.Range("M107").Value = 1194.73334
$ws.Range("N107").Value = -4739
$ws.Range("H129").Value = 80780
$ws.Range("J129").Value = 80780
$ws.Range("L129").Value = 80780
$ws.Range("N129").Value = -90780
$ws.Range("H130").Value = 20709
$ws.Range("I130").Value = 20709
$ws.Range("K130").Value = 20709
$ws.Range("M130").Value = -15689
$ws.Range("H134").Value = 1543.7826
$ws.Range("I134").Value = 1340.591
$ws.Range("K134").Value = 4021.773
$ws.Range("M134").Value = -1486.773

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 882.3333
$ws.Range("I16").Value = 867.625
$ws.Range("K16").Value = 867.625
$ws.Range("M16").Value = -580.625
$ws.Range("H62").Value = 68949
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 82038.8
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 82038.8
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -83286.8
$ws.Range("H65").Value = 68949
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 82038.8
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 410194
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -416434
$ws.Range("H99").Value = 14294.448
$ws.Range("I99").Value = 10698.5
$ws.Range("K99").Value = 10698.5
$ws.Range("M99").Value = -9200.5
$ws.Range("H113").Value = 882.3333
$ws.Range("I113").Value = 867.625
$ws.Range("K113").Value = 867.625
$ws.Range("M113").Value = 1302.375
$ws.Range("H126").Value = 14294.448
$ws.Range("I126").Value = 10698.5
$ws.Range("K126").Value = 32095.5
$ws.Range("M126").Value = -29625.5
$ws.Range("H132").Value = 3247.8572
$ws.Range("I132").Value = 1750.8
$ws.Range("K132").Value = 5252.4
$ws.Range("M132").Value = -2722.4
$ws.Range("H134").Value = 3015.4736
$ws.Range("I134").Value = 2485.3333
$ws.Range("J134").Value = 5003.5
$ws.Range("K134").Value = 7455.999899999999
$ws.Range("L134").Value = 15010.5
$ws.Range("M134").Value = -4920.999899999999
$ws.Range("N134").Value = -20080.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71498.14
$ws.Range("I2").Value = 90943.09
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 545658.54
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -545545.54
$ws.Range("N2").Value = -1426
$ws.Range("H36").Value = 2000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H107").Value = 569.51514
$ws.Range("I107").Value = 133
$ws.Range("J107").Value = 613.1667
$ws.Range("K107").Value = 399
$ws.Range("L107").Value = 1839.5001
$ws.Range("M107").Value = 1521
$ws.Range("N107").Value = -5679.5001
$ws.Range("H113").Value = 989.2727
$ws.Range("J113").Value = 485.375
$ws.Range("L113").Value = 1456.125
$ws.Range("N113").Value = -5796.125

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1352.5
$ws.Range("J107").Value = 946.6667
$ws.Range("L107").Value = 946.6667
$ws.Range("N107").Value = -4786.6667
$ws.Range("H122").Value = 55478.473
$ws.Range("I122").Value = 2318.5625
$ws.Range("K122").Value = 6955.6875
$ws.Range("M122").Value = -4505.6875
$ws.Range("H132").Value = 3378.625
$ws.Range("I132").Value = 751.5
$ws.Range("K132").Value = 2254.5
$ws.Range("M132").Value = 275.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4011.2
$ws.Range("I16").Value = 2014
$ws.Range("K16").Value = 2014
$ws.Range("M16").Value = -1844
$ws.Range("H55").Value = 846
$ws.Range("I55").Value = 811.375
$ws.Range("K55").Value = 811.375
$ws.Range("M55").Value = -638.375
$ws.Range("H136").Value = 2450
$ws.Range("I136").Value = 2450
$ws.Range("K136").Value = 7350
$ws.Range("M136").Value = -4800

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H60").Value = 39000
$ws.Range("J60").Value = 39000
$ws.Range("L60").Value = 39000
$ws.Range("N60").Value = -40644
$ws.Range("H113").Value = 646.5
$ws.Range("I113").Value = 780.6667
$ws.Range("J113").Value = 244
$ws.Range("K113").Value = 2342.0001
$ws.Range("L113").Value = 732
$ws.Range("M113").Value = -172.0001000000002
$ws.Range("N113").Value = -5072
$ws.Range("H126").Value = 7082.3335
$ws.Range("I126").Value = 5246
$ws.Range("K126").Value = 15738
$ws.Range("M126").Value = -13268
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3593.6667
$ws.Range("I136").Value = 1410.5
$ws.Range("K136").Value = 4231.5
$ws.Range("M136").Value = -1681.5
